$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(301, 6, 45, 30, 60, 45)
    3  = @(501, 9, 52, 30, 75, 45)
    4  = @(701, 3, 90, 45, 97, 15)
    5  = @(201, 9, 30, 15, 45, 30)
    6  = @(1202, 2, 10, 10, 10, 10)
    7  = @(1203, 3, 15, 15, 15, 15)
    8  = @(101, 9, 30, 15, 60, 15)
    9  = @(901, 16, 15, 45, 60, 60)
    10 = @(1201, 2, 10, 10, 10, 10)
    11 = @(1001, 18, 30, 75, 60, 72)
    12 = @(401, 9, 48, 67, 75, 45)
    13 = @(601, 9, 60, 67, 60, 42)
    14 = @(801, 3, 67, 65, 52, 45)
    15 = @(902, 1, 0, 0, 0, 0)
    16 = @(1, 0, 2, 2, 2, 2)
    17 = @(2, 0, 2, 2, 2, 2)
    18 = @(1101, 0, 15, 30, 30, 0)
    20 = @(3, 0, 3, 3, 3, 3)
    21 = @(502, 0, 4, 0, 0, 0)
}

foreach ($r in $data.Keys) {
    $values = $data[$r]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $values[$i]
    }
}
